$d = $word.ActiveDocument

# --- Reformat table captions: add a dedicated "TableCaption" paragraph
# style (based on the document's Normal1 style) so table captions get a
# bold, smaller, "captiony" look. ------------------------------------
$style = $d.Styles.Add("TableCaption", 1)
$style.BaseStyle = "Normal1"
$style.QuickStyle = $true

$style.Font.NameAscii = "Calibri"
$style.Font.Name = "Calibri"
$style.Font.NameFarEast = "Calibri"
$style.Font.NameBi = "Calibri"
$style.Font.Bold = $true
$style.Font.Color = 3355443
$style.Font.Size = 9

$style.ParagraphFormat.LineSpacingRule = 5
$style.ParagraphFormat.LineSpacing = 15.6

# --- Word tracks the location of the user's last edit with the special
# "_GoBack" bookmark; saving after this edit leaves it sitting right
# before the {{table:ecosystems}} placeholder (and removes the old one
# that used to sit around {{table:protection}}). ----------------------
$r = $d.Content
$r.Find.Execute("{{table:ecosystems}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r)
